$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.1565217391304348
$ws.Cells.Item(2, 3).Value = 0.591304347826087
$ws.Cells.Item(2, 10).Value = 0.008695652173913044
$ws.Cells.Item(2, 16).Value = 0.09565217391304348
$ws.Cells.Item(2, 19).Value = 0.1478260869565217

# Row 3
$ws.Cells.Item(3, 2).Value = 0.01408450704225352
$ws.Cells.Item(3, 3).Value = 0.04225352112676056
$ws.Cells.Item(3, 10).Value = 0.01408450704225352
$ws.Cells.Item(3, 16).Value = 0.676056338028169
$ws.Cells.Item(3, 19).Value = 0.2535211267605634

# Row 4
$ws.Cells.Item(4, 16).Value = 0.68
$ws.Cells.Item(4, 19).Value = 0.32

# Row 6
$ws.Cells.Item(6, 2).Value = 0.07547169811320754
$ws.Cells.Item(6, 6).Value = 0.1320754716981132
$ws.Cells.Item(6, 10).Value = 0.2452830188679245
$ws.Cells.Item(6, 15).Value = 0.01886792452830189
$ws.Cells.Item(6, 17).Value = 0.1037735849056604
$ws.Cells.Item(6, 18).Value = 0.05660377358490566
$ws.Cells.Item(6, 19).Value = 0.3679245283018868

# Row 7
$ws.Cells.Item(7, 2).Value = 0.03488372093023256
$ws.Cells.Item(7, 4).Value = 0.02325581395348837
$ws.Cells.Item(7, 6).Value = 0.05813953488372093
$ws.Cells.Item(7, 10).Value = 0.1395348837209302
$ws.Cells.Item(7, 15).Value = 0.02325581395348837
$ws.Cells.Item(7, 17).Value = 0.1744186046511628
$ws.Cells.Item(7, 18).Value = 0.1162790697674419
$ws.Cells.Item(7, 19).Value = 0.4302325581395349

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1368421052631579
$ws.Cells.Item(8, 4).Value = 0.005263157894736842
$ws.Cells.Item(8, 6).Value = 0.06842105263157895
$ws.Cells.Item(8, 10).Value = 0.1210526315789474
$ws.Cells.Item(8, 15).Value = 0.03684210526315789
$ws.Cells.Item(8, 17).Value = 0.1578947368421053
$ws.Cells.Item(8, 18).Value = 0.06842105263157895
$ws.Cells.Item(8, 19).Value = 0.4052631578947368

# Row 9
$ws.Cells.Item(9, 2).Value = 0.08620689655172414
$ws.Cells.Item(9, 4).Value = 0.05172413793103448
$ws.Cells.Item(9, 6).Value = 0.03448275862068965
$ws.Cells.Item(9, 10).Value = 0.08620689655172414
$ws.Cells.Item(9, 15).Value = 0.03448275862068965
$ws.Cells.Item(9, 17).Value = 0.1206896551724138
$ws.Cells.Item(9, 18).Value = 0.103448275862069
$ws.Cells.Item(9, 19).Value = 0.4827586206896552

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1055776892430279
$ws.Cells.Item(10, 4).Value = 0.03784860557768924
$ws.Cells.Item(10, 5).Value = 0.00398406374501992
$ws.Cells.Item(10, 6).Value = 0.07171314741035857
$ws.Cells.Item(10, 10).Value = 0.1215139442231076
$ws.Cells.Item(10, 15).Value = 0.0298804780876494
$ws.Cells.Item(10, 17).Value = 0.149402390438247
$ws.Cells.Item(10, 18).Value = 0.08366533864541832
$ws.Cells.Item(10, 19).Value = 0.3964143426294821

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1203007518796992
$ws.Cells.Item(11, 10).Value = 0.06766917293233082
$ws.Cells.Item(11, 11).Value = 0.1654135338345865
$ws.Cells.Item(11, 12).Value = 0.6165413533834586
$ws.Cells.Item(11, 19).Value = 0.03007518796992481

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7317073170731707
$ws.Cells.Item(12, 10).Value = 0.2073170731707317
$ws.Cells.Item(12, 11).Value = 0.01219512195121951
$ws.Cells.Item(12, 12).Value = 0.02439024390243903
$ws.Cells.Item(12, 19).Value = 0.02439024390243903

# Row 13
$ws.Cells.Item(13, 7).Value = 0.7368421052631579
$ws.Cells.Item(13, 10).Value = 0.2631578947368421

# Row 14
$ws.Cells.Item(14, 7).Value = 1

# Row 15
$ws.Cells.Item(15, 6).Value = 0.008771929824561403
$ws.Cells.Item(15, 8).Value = 0.1666666666666667
$ws.Cells.Item(15, 9).Value = 0.03508771929824561
$ws.Cells.Item(15, 10).Value = 0.3333333333333333
$ws.Cells.Item(15, 11).Value = 0.03508771929824561
$ws.Cells.Item(15, 13).Value = 0.02631578947368421
$ws.Cells.Item(15, 15).Value = 0.1140350877192982
$ws.Cells.Item(15, 19).Value = 0.2807017543859649

# Row 16
$ws.Cells.Item(16, 6).Value = 0.02666666666666667
$ws.Cells.Item(16, 8).Value = 0.2533333333333334
$ws.Cells.Item(16, 9).Value = 0.06666666666666667
$ws.Cells.Item(16, 10).Value = 0.32
$ws.Cells.Item(16, 11).Value = 0.1066666666666667
$ws.Cells.Item(16, 15).Value = 0.05333333333333334
$ws.Cells.Item(16, 19).Value = 0.1733333333333333

# Row 17
$ws.Cells.Item(17, 6).Value = 0.02877697841726619
$ws.Cells.Item(17, 8).Value = 0.2446043165467626
$ws.Cells.Item(17, 9).Value = 0.04316546762589928
$ws.Cells.Item(17, 10).Value = 0.3381294964028777
$ws.Cells.Item(17, 11).Value = 0.08633093525179857
$ws.Cells.Item(17, 13).Value = 0.02877697841726619
$ws.Cells.Item(17, 15).Value = 0.07194244604316546
$ws.Cells.Item(17, 19).Value = 0.158273381294964

# Row 18
$ws.Cells.Item(18, 6).Value = 0.05194805194805195
$ws.Cells.Item(18, 8).Value = 0.1688311688311688
$ws.Cells.Item(18, 9).Value = 0.07792207792207792
$ws.Cells.Item(18, 10).Value = 0.4415584415584415
$ws.Cells.Item(18, 11).Value = 0.1168831168831169
$ws.Cells.Item(18, 15).Value = 0.05194805194805195
$ws.Cells.Item(18, 19).Value = 0.09090909090909091

# Row 19
$ws.Cells.Item(19, 6).Value = 0.02554744525547445
$ws.Cells.Item(19, 8).Value = 0.1934306569343066
$ws.Cells.Item(19, 9).Value = 0.06569343065693431
$ws.Cells.Item(19, 10).Value = 0.3759124087591241
$ws.Cells.Item(19, 11).Value = 0.1441605839416058
$ws.Cells.Item(19, 13).Value = 0.02372262773722628
$ws.Cells.Item(19, 14).Value = 0.001824817518248175
$ws.Cells.Item(19, 15).Value = 0.08029197080291971
$ws.Cells.Item(19, 19).Value = 0.08941605839416059
